# Consolidate the multiple text runs in the caption textbox ("Followed" /
# " " / "by" / " " / "a" / " " / "picture") into a single run reading
# "Followed by a picture".
#
# Note: the caption's runs already concatenate to exactly that string, so
# assigning the identical text back to TextRange.Text is treated as a
# no-op by the text engine (nothing to reconcile, so no run-merge happens).
# To force a genuine rewrite of the run structure we first set the text to
# a different placeholder value (collapsing everything into one run) and
# then set it to the desired final text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 3") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(3)
}

$shape.TextFrame.TextRange.Text = "."
$shape.TextFrame.TextRange.Text = "Followed by a picture"
